$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force the Price/Volume columns to Text format so that
# numeric-looking strings (e.g. "51.759.40", "1.00") are stored as text,
# matching the original inline-string cell contents, not auto-converted
# to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$data = @(
    @(2, "Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "51.882.62", "  +0.30%  "),
    @(3, "Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "2.816.00", "  +1.49%  "),
    @(4, "TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "0.999", "  -0.11%  "),
    @(5, "BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "356.82", "  -0.04%  "),
    @(6, "Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "109.75", "  +0.44%  "),
    @(7, "XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.556", "  -0.10%  "),
    @(8, "USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "1.00", "  +0.03%  "),
    @(9, "Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.634", "  +8.04%  "),
    @(10, "Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "40.09", "  +0.58%  "),
    @(11, "TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.135", "  +0.03%  "),
    @(12, "Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.0838", "  -0.72%  "),
    @(13, "Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "20.01", "  +2.75%  "),
    @(14, "Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "7.81", "  +2.63%  "),
    @(15, "WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "3.246.18", "  +1.05%  "),
    @(16, "WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "2.805.49", "  +1.57%  "),
    @(17, "Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "0.945", "  +1.04%  "),
    @(18, "WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "51.853.04", "  +0.30%  "),
    @(19, "Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "7.72", "  +3.73%  "),
    @(20, "ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "3.17", "  +3.40%  "),
    @(21, "InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "13.68", "  +4.22%  "),
    @(22, "ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.0₃0980", "  +0.95%  "),
    @(23, "Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "70.52", "  +0.41%  "),
    @(24, "BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "268.98", "  +0.06%  "),
    @(25, "PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "2.77", "  +0.96%  "),
    @(26, "EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "26.22", "  -0.57%  "),
    @(27, "Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "1.00", "  +0.12%  "),
    @(28, "Kaspa", "https://coinranking.com/coin/V8GxkwWow+kaspa-kas", "0.164", "  +0.47%  "),
    @(29, "Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "10.41", "  +1.68%  "),
    @(30, "InjectiveProtocol", "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj", "37.94", "  +9.47%  "),
    @(31, "Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "6.21", "  +0.68%  "),
    @(32, "OKB", "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb", "52.00", "  +0.34%  "),
    @(33, "Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "2.11", "  -4.85%  "),
    @(34, "RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "5.69", "  +11.06%  "),
    @(35, "VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.0446", "  -1.03%  "),
    @(36, "Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.0865", "  +3.10%  "),
    @(37, "FirstDigitalUSD", "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd", "0.999", "  -0.16%  "),
    @(38, "Celestia", "https://coinranking.com/coin/YQcD0lBl7+celestia-tia", "18.94", "  +1.31%  "),
    @(39, "ARBITRUM", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", "2.01", "  +2.64%  "),
    @(40, "LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "3.16", "  +0.82%  "),
    @(41, "Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.115", "  +1.27%  "),
    @(42, "Stacks", "https://coinranking.com/coin/mMPrMcB7+stacks-stx", "2.50", "  -1.36%  "),
    @(43, "EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "22.09", "  +1.91%  "),
    @(44, "WEMIXToken", "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix", "2.20", "  -0.98%  "),
    @(45, "Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "119.20", "  -0.52%  "),
    @(46, "ApeXProtocol", "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex", "2.47", "  +8.52%  "),
    @(47, "NEARProtocol", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near", "3.40", "  +4.15%  "),
    @(48, "Maker", "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr", "2.111.38", "  +1.31%  "),
    @(49, "SEI", "https://coinranking.com/coin/8nxCqs-uj+sei-sei", "0.927", "  -0.90%  "),
    @(50, "TrustWalletToken", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt", "1.37", "  +9.22%  "),
    @(51, "THORChain", "https://coinranking.com/coin/ybmU-kKU+thorchain-rune", "5.45", "  -5.18%  "),
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 2).Value = $item[1]
    $ws.Cells.Item($r, 3).Value = $item[2]
    $ws.Cells.Item($r, 4).Value = $item[3]
    $ws.Cells.Item($r, 5).Value = $item[4]
}

# Restore the original (unstyled) formatting for the data cells now that
# the text values are safely in place.
$ws.Range("D2:E51").ClearFormats()
